$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("20160405_164902", 2331.295, 'convert unicode to ascii, remove multiple spaces, convert to lower, trim "space" and ","', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.992, 0.99009900990099, '0 filters: ', 0.354166666666667),
    @("20160405_172754", 2423.942, 'convert unicode to ascii, remove multiple spaces, convert to lower, trim "space" and ","', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.992, 0.99009900990099, '0 filters: ', 0.395833333333333),
    @("20160405_180818", 2372.816, 'convert unicode to ascii, remove multiple spaces, convert to lower, trim "space" and ","', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.994666666666667, 0.99009900990099, '0 filters: ', 0.40625),
    @("20160405_184750", 2433.274, 'convert unicode to ascii, remove multiple spaces, convert to lower, trim "space" and ","', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.992666666666667, 0.99009900990099, '0 filters: ', 0.40625),
    @("20160405_192824", 2391.304, 'convert unicode to ascii, remove multiple spaces, convert to lower, trim "space" and ","', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.993333333333333, 0.99009900990099, '0 filters: ', 0.375),
    @("20160406_081448", 3327.672, 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.992666666666667, 0.99009900990099, '0 filters: ', 0.385416666666667),
    @("20160406_091015", 3384.742, 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.996, 0.99009900990099, '0 filters: ', 0.395833333333333),
    @("20160406_100640", 4688.205, 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.994, 0.99009900990099, '0 filters: ', 0.395833333333333),
    @("20160406_112448", 5738.497, 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.992, 0.99009900990099, '0 filters: ', 0.385416666666667),
    @("20160406_130027", 6565.086, 'convert to lower, convert unicode to ascii, trim "space" and ",", remove multiple spaces', '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, #(, #+, #/', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.989333333333333, 0.99009900990099, '0 filters: ', 0.385416666666667)
)

$startRow = 22
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
